# Generate Report for Handoff
# - Flip the per-language "Status" cells from "Handed back: in sync with en-US"
#   to "Ready for handoff" on all three sheets.
# - Bump the associated handoff/generate timestamps.
# - Narrow the Status/date columns (they no longer need to fit the old, longer
#   "Handed back: in sync with en-US" text).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------------
$ovw = $wb.Sheets.Item("Overview")
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$ovw.Range("G2").Value = "2016-08-27 16:57:55"

$ovw.Columns.Item(5).ColumnWidth = 16.33
$ovw.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-27 16:57:51"

$zhcn.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet ------------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-27 16:57:55"

$dede.Columns.Item(3).ColumnWidth = 16.33
